$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.040094350233371
$ws.Cells.Item(2, 4).Value = 1.050244674905319
$ws.Cells.Item(2, 5).Value = 1.048547495786166
$ws.Cells.Item(2, 6).Value = 1.0594430891009
$ws.Cells.Item(2, 9).Value = 1.036478988838291
$ws.Cells.Item(2, 10).Value = 1.045182815016401
$ws.Cells.Item(2, 11).Value = 1.05299948975008
$ws.Cells.Item(2, 12).Value = 1.051307035798264
$ws.Cells.Item(2, 13).Value = 1.06217257962361
$ws.Cells.Item(2, 14).Value = 1.018962365659191

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.040982097413727
$ws.Cells.Item(3, 4).Value = 1.051063406903668
$ws.Cells.Item(3, 5).Value = 1.049333520566056
$ws.Cells.Item(3, 6).Value = 1.060303557600974
$ws.Cells.Item(3, 9).Value = 1.036584975585148
$ws.Cells.Item(3, 10).Value = 1.045716346847191
$ws.Cells.Item(3, 11).Value = 1.053630613117952
$ws.Cells.Item(3, 12).Value = 1.051905195295929
$ws.Cells.Item(3, 13).Value = 1.062847161545223
$ws.Cells.Item(3, 14).Value = 1.019141258202088

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.041557170494382
$ws.Cells.Item(4, 4).Value = 1.051594096964151
$ws.Cells.Item(4, 5).Value = 1.049843068205098
$ws.Cells.Item(4, 6).Value = 1.06086133609219
$ws.Cells.Item(4, 9).Value = 1.036652280463742
$ws.Cells.Item(4, 10).Value = 1.046061567461547
$ws.Cells.Item(4, 11).Value = 1.054039255072296
$ws.Cells.Item(4, 12).Value = 1.052292522002801
$ws.Cells.Item(4, 13).Value = 1.063284013216996
$ws.Cells.Item(4, 14).Value = 1.019256959698746

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.041799083238834
$ws.Cells.Item(5, 4).Value = 1.051817416503504
$ws.Cells.Item(5, 5).Value = 1.050057504598612
$ws.Cells.Item(5, 6).Value = 1.061096063196148
$ws.Cells.Item(5, 9).Value = 1.036680269420405
$ws.Cells.Item(5, 10).Value = 1.046206694505678
$ws.Cells.Item(5, 11).Value = 1.054211109635432
$ws.Cells.Item(5, 12).Value = 1.052455419183975
$ws.Cells.Item(5, 13).Value = 1.063467748378264
$ws.Cells.Item(5, 14).Value = 1.019305587218885

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.041839710324009
$ws.Cells.Item(6, 4).Value = 1.051854925535579
$ws.Cells.Item(6, 5).Value = 1.050093522412708
$ws.Cells.Item(6, 6).Value = 1.061135488755091
$ws.Cells.Item(6, 9).Value = 1.036684950926408
$ws.Cells.Item(6, 10).Value = 1.046231061720947
$ws.Cells.Item(6, 11).Value = 1.05423996833584
$ws.Cells.Item(6, 12).Value = 1.052482774110262
$ws.Cells.Item(6, 13).Value = 1.063498603125451
$ws.Cells.Item(6, 14).Value = 1.019313751194517

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.041560402348853
$ws.Cells.Item(7, 4).Value = 1.051597080118625
$ws.Cells.Item(7, 5).Value = 1.049845932642368
$ws.Cells.Item(7, 6).Value = 1.060864471598223
$ws.Cells.Item(7, 9).Value = 1.036652655656757
$ws.Cells.Item(7, 10).Value = 1.046063506672355
$ws.Cells.Item(7, 11).Value = 1.054041551161618
$ws.Cells.Item(7, 12).Value = 1.052294698389134
$ws.Cells.Item(7, 13).Value = 1.063286467972153
$ws.Cells.Item(7, 14).Value = 1.019257609515257

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.040394234859222
$ws.Cells.Item(8, 4).Value = 1.050521178540368
$ws.Cells.Item(8, 5).Value = 1.048812941458176
$ws.Cells.Item(8, 6).Value = 1.059733680734506
$ws.Cells.Item(8, 9).Value = 1.036515071308028
$ws.Cells.Item(8, 10).Value = 1.045363125935495
$ws.Cells.Item(8, 11).Value = 1.05321272552649
$ws.Cells.Item(8, 12).Value = 1.051509128295934
$ws.Cells.Item(8, 13).Value = 1.062400483904667
$ws.Cells.Item(8, 14).Value = 1.019022834114624

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.038344276922276
$ws.Cells.Item(9, 4).Value = 1.048632394660416
$ws.Cells.Item(9, 5).Value = 1.046999933676211
$ws.Cells.Item(9, 6).Value = 1.057748807086954
$ws.Cells.Item(9, 9).Value = 1.036262887922536
$ws.Cells.Item(9, 10).Value = 1.044128943980085
$ws.Cells.Item(9, 11).Value = 1.051754307216388
$ws.Cells.Item(9, 12).Value = 1.050127040905666
$ws.Cells.Item(9, 13).Value = 1.060842032355279
$ws.Cells.Item(9, 14).Value = 1.018608737747693

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.036981074731623
$ws.Cells.Item(10, 4).Value = 1.047378080076904
$ws.Cells.Item(10, 5).Value = 1.045796244572561
$ws.Cells.Item(10, 6).Value = 1.056430862010294
$ws.Cells.Item(10, 9).Value = 1.036088252789818
$ws.Cells.Item(10, 10).Value = 1.043306217064122
$ws.Cells.Item(10, 11).Value = 1.050783512674205
$ws.Cells.Item(10, 12).Value = 1.049207203609906
$ws.Cells.Item(10, 13).Value = 1.059805017274768
$ws.Cells.Item(10, 14).Value = 1.01833243814576

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.036391626910109
$ws.Cells.Item(11, 4).Value = 1.046836126615912
$ws.Cells.Item(11, 5).Value = 1.045276237907695
$ws.Cells.Item(11, 6).Value = 1.055861458592385
$ws.Cells.Item(11, 9).Value = 1.036011096895882
$ws.Cells.Item(11, 10).Value = 1.042949997597145
$ws.Cells.Item(11, 11).Value = 1.050363518281281
$ws.Cells.Item(11, 12).Value = 1.04880929015264
$ws.Cells.Item(11, 13).Value = 1.059356461052096
$ws.Cells.Item(11, 14).Value = 1.018212747493055

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.036172805587781
$ws.Cells.Item(12, 4).Value = 1.046634999012329
$ws.Cells.Item(12, 5).Value = 1.04508326615618
$ws.Cells.Item(12, 6).Value = 1.055650150321235
$ws.Cells.Item(12, 9).Value = 1.035982207344639
$ws.Cells.Item(12, 10).Value = 1.042817687111182
$ws.Cells.Item(12, 11).Value = 1.050207570127072
$ws.Cells.Item(12, 12).Value = 1.048661546206066
$ws.Cells.Item(12, 13).Value = 1.059189920712636
$ws.Cells.Item(12, 14).Value = 1.018168281854225

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.036219737785604
$ws.Cells.Item(13, 4).Value = 1.046678133480233
$ws.Cells.Item(13, 5).Value = 1.045124650992855
$ws.Cells.Item(13, 6).Value = 1.055695467875709
$ws.Cells.Item(13, 9).Value = 1.035988414668501
$ws.Cells.Item(13, 10).Value = 1.042846067900006
$ws.Cells.Item(13, 11).Value = 1.050241018956851
$ws.Cells.Item(13, 12).Value = 1.048693235099333
$ws.Cells.Item(13, 13).Value = 1.059225640833346
$ws.Cells.Item(13, 14).Value = 1.018177820202922

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.036373536495551
$ws.Cells.Item(14, 4).Value = 1.046819497702996
$ws.Cells.Item(14, 5).Value = 1.04526028308192
$ws.Cells.Item(14, 6).Value = 1.055843987826327
$ws.Cells.Item(14, 9).Value = 1.036008713573469
$ws.Cells.Item(14, 10).Value = 1.042939060655943
$ws.Cells.Item(14, 11).Value = 1.050350626397872
$ws.Cells.Item(14, 12).Value = 1.048797076383671
$ws.Cells.Item(14, 13).Value = 1.059342693268949
$ws.Cells.Item(14, 14).Value = 1.018209072093628

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.03646831367684
$ws.Cells.Item(15, 4).Value = 1.046906620512647
$ws.Cells.Item(15, 5).Value = 1.045343874642484
$ws.Cells.Item(15, 6).Value = 1.055935521568047
$ws.Cells.Item(15, 9).Value = 1.036021189884549
$ws.Cells.Item(15, 10).Value = 1.042996357300326
$ws.Cells.Item(15, 11).Value = 1.050418166685693
$ws.Cells.Item(15, 12).Value = 1.048861064263768
$ws.Cells.Item(15, 13).Value = 1.059414822912928
$ws.Cells.Item(15, 14).Value = 1.018228326474947

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.037020212008303
$ws.Cells.Item(16, 4).Value = 1.047414072611952
$ws.Cells.Item(16, 5).Value = 1.045830781109681
$ws.Cells.Item(16, 6).Value = 1.056468678473758
$ws.Cells.Item(16, 9).Value = 1.036093341026213
$ws.Cells.Item(16, 10).Value = 1.043329858848325
$ws.Cells.Item(16, 11).Value = 1.050811394156461
$ws.Cells.Item(16, 12).Value = 1.049233619987896
$ws.Cells.Item(16, 13).Value = 1.059834796682663
$ws.Cells.Item(16, 14).Value = 1.018340380584194

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.037366625904559
$ws.Cells.Item(17, 4).Value = 1.047732699237044
$ws.Cells.Item(17, 5).Value = 1.046136526853342
$ws.Cells.Item(17, 6).Value = 1.056803456530314
$ws.Cells.Item(17, 9).Value = 1.036138188267414
$ws.Cells.Item(17, 10).Value = 1.043539063608257
$ws.Cells.Item(17, 11).Value = 1.051058154398705
$ws.Cells.Item(17, 12).Value = 1.049467417663988
$ws.Cells.Item(17, 13).Value = 1.060098364260496
$ws.Cells.Item(17, 14).Value = 1.018410655821418

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.037568763032641
$ws.Cells.Item(18, 4).Value = 1.047918661815987
$ws.Cells.Item(18, 5).Value = 1.04631497870601
$ws.Cells.Item(18, 6).Value = 1.056998849847512
$ws.Cells.Item(18, 9).Value = 1.036164198494964
$ws.Cells.Item(18, 10).Value = 1.043661091697795
$ws.Cells.Item(18, 11).Value = 1.051202120693177
$ws.Cells.Item(18, 12).Value = 1.049603824671221
$ws.Cells.Item(18, 13).Value = 1.060252144798262
$ws.Cells.Item(18, 14).Value = 1.018451641215825

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.037637700053655
$ws.Cells.Item(19, 4).Value = 1.047982089371515
$ws.Cells.Item(19, 5).Value = 1.046375845718065
$ws.Cells.Item(19, 6).Value = 1.057065494714901
$ws.Cells.Item(19, 9).Value = 1.036173042119041
$ws.Cells.Item(19, 10).Value = 1.04370270050696
$ws.Cells.Item(19, 11).Value = 1.051251215408417
$ws.Cells.Item(19, 12).Value = 1.049650342130864
$ws.Cells.Item(19, 13).Value = 1.060304587746911
$ws.Cells.Item(19, 14).Value = 1.01846561531994

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.037329450710885
$ws.Cells.Item(20, 4).Value = 1.047698501908325
$ws.Cells.Item(20, 5).Value = 1.04610371126552
$ws.Cells.Item(20, 6).Value = 1.056767525283225
$ws.Cells.Item(20, 9).Value = 1.036133391929154
$ws.Cells.Item(20, 10).Value = 1.04351661765042
$ws.Cells.Item(20, 11).Value = 1.051031675706453
$ws.Cells.Item(20, 12).Value = 1.049442329571628
$ws.Cells.Item(20, 13).Value = 1.060070081181412
$ws.Cells.Item(20, 14).Value = 1.018403116463398

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.036328243122932
$ws.Cells.Item(21, 4).Value = 1.046777864533449
$ws.Cells.Item(21, 5).Value = 1.045220337777147
$ws.Cells.Item(21, 6).Value = 1.055800247045671
$ws.Cells.Item(21, 9).Value = 1.036002742407627
$ws.Cells.Item(21, 10).Value = 1.042911676451011
$ws.Cells.Item(21, 11).Value = 1.050318348167919
$ws.Cells.Item(21, 12).Value = 1.048766496078307
$ws.Cells.Item(21, 13).Value = 1.05930822220488
$ws.Cells.Item(21, 14).Value = 1.018199869385713

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.03569947300717
$ws.Cells.Item(22, 4).Value = 1.04620005426831
$ws.Cells.Item(22, 5).Value = 1.044665978975241
$ws.Cells.Item(22, 6).Value = 1.05519320135266
$ws.Cells.Item(22, 9).Value = 1.0359192650495
$ws.Cells.Item(22, 10).Value = 1.042531357114179
$ws.Cells.Item(22, 11).Value = 1.049870178240517
$ws.Cells.Item(22, 12).Value = 1.048341913800557
$ws.Cells.Item(22, 13).Value = 1.058829636757351
$ws.Cells.Item(22, 14).Value = 1.018072038452795

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.036032726359479
$ws.Cells.Item(23, 4).Value = 1.046506264092709
$ws.Cells.Item(23, 5).Value = 1.044959754638369
$ws.Cells.Item(23, 6).Value = 1.055514900918938
$ws.Cells.Item(23, 9).Value = 1.035963644108477
$ws.Cells.Item(23, 10).Value = 1.042732968208153
$ws.Cells.Item(23, 11).Value = 1.05010773007374
$ws.Cells.Item(23, 12).Value = 1.048566960098686
$ws.Cells.Item(23, 13).Value = 1.059083303096414
$ws.Cells.Item(23, 14).Value = 1.018139807816948

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.037346248334424
$ws.Cells.Item(24, 4).Value = 1.047713953860052
$ws.Cells.Item(24, 5).Value = 1.046118538860369
$ws.Cells.Item(24, 6).Value = 1.056783760687504
$ws.Cells.Item(24, 9).Value = 1.036135559646589
$ws.Cells.Item(24, 10).Value = 1.043526760005038
$ws.Cells.Item(24, 11).Value = 1.051043640179069
$ws.Cells.Item(24, 12).Value = 1.049453665687191
$ws.Cells.Item(24, 13).Value = 1.060082860945385
$ws.Cells.Item(24, 14).Value = 1.01840652319015

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.038873641050199
$ws.Cells.Item(25, 4).Value = 1.049119839155149
$ws.Cells.Item(25, 5).Value = 1.047467769095145
$ws.Cells.Item(25, 6).Value = 1.058261017730661
$ws.Cells.Item(25, 9).Value = 1.036329234040062
$ws.Cells.Item(25, 10).Value = 1.04444800437293
$ws.Cells.Item(25, 11).Value = 1.052131088051812
$ws.Cells.Item(25, 12).Value = 1.050484075786354
$ws.Cells.Item(25, 13).Value = 1.061244592098414
$ws.Cells.Item(25, 14).Value = 1.018715835137052
